$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.651.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.58%  "

$ws.Range("D3").Value = "'2.251.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.95%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'234.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").Value = "'0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "

$ws.Range("D7").Value = "'63.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.52%  "

$ws.Range("D10").Value = "'59.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.65%  "

$ws.Range("D11").Value = "'0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.89%  "

$ws.Range("D12").Value = "'0.105"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").Value = "'2.584.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.87%  "

$ws.Range("D14").Value = "'16.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "'22.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.82%  "

$ws.Range("D16").Value = "'0.823"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").Value = "'2.256.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("D19").Value = "'41.480.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.03%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'74.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.51%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +9.53%  "

$ws.Range("D22").Value = "'6.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("D23").Value = "'251.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.48%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.13%  "

$ws.Range("D26").Value = "'2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("D27").Value = "'9.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.32%  "

$ws.Range("D28").Value = "'0.148"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.50%  "

$ws.Range("D29").Value = "'170.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "'20.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.97%  "

$ws.Range("D31").Value = "'1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").Value = "'2.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.02%  "

$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("E34").Value = "  +8.54%  "

$ws.Range("D35").Value = "'4.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.68%  "

$ws.Range("D36").Value = "'0.0637"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "

$ws.Range("D37").Value = "'6.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.40%  "

$ws.Range("D38").Value = "'3.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.93%  "

$ws.Range("D39").Value = "'2.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").Value = "'0.000262"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +60.73%  "

$ws.Range("D41").Value = "'5.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.80%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  +5.59%  "

$ws.Range("D44").Value = "'8.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.77%  "

$ws.Range("D45").Value = "'102.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").Value = "'0.0994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.28%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'17.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("D49").Value = "'1.501.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("E51").Value = "  -0.65%  "
